$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
  $ws.Cells.Item(2, 1).Value2 = "FAPs"
  $ws.Cells.Item(2, 2).Value2 = "Rspo3"
  $ws.Cells.Item(2, 3).Value2 = "Lgr4"
  $ws.Cells.Item(2, 4).Value2 = "ECs"
  $ws.Cells.Item(2, 5).Value2 = 3
  $ws.Cells.Item(2, 6).Value2 = 1
  $ws.Cells.Item(2, 7).Value2 = 4.83236
  $ws.Cells.Item(2, 8).Value2 = 14.49708
  $ws.Cells.Item(2, 9).Value2 = 0.975350813525687
  $ws.Cells.Item(2, 10).Value2 = 0.975350813525687
  $ws.Cells.Item(2, 11).Value2 = 3
  $ws.Cells.Item(2, 12).Value2 = 1
  $ws.Cells.Item(2, 13).Value2 = 2.451301666666667
  $ws.Cells.Item(2, 14).Value2 = 7.353905
  $ws.Cells.Item(2, 15).Value2 = 0.1099648918267336
  $ws.Cells.Item(2, 16).Value2 = 0.1099648918267337
  $ws.Cells.Item(2, 17).Value2 = 11.84557212193334
  $ws.Cells.Item(2, 18).Value2 = 106.6101490974
  $ws.Cells.Item(2, 19).Value2 = 0.1072543467024688
  $ws.Cells.Item(2, 20).Value2 = 0.1072543467024688

# Row 3
  $ws.Cells.Item(3, 1).Value2 = "FAPs"
  $ws.Cells.Item(3, 2).Value2 = "Rspo3"
  $ws.Cells.Item(3, 3).Value2 = "Lgr4"
  $ws.Cells.Item(3, 4).Value2 = "FAPs"
  $ws.Cells.Item(3, 5).Value2 = 3
  $ws.Cells.Item(3, 6).Value2 = 1
  $ws.Cells.Item(3, 7).Value2 = 4.83236
  $ws.Cells.Item(3, 8).Value2 = 14.49708
  $ws.Cells.Item(3, 9).Value2 = 0.975350813525687
  $ws.Cells.Item(3, 10).Value2 = 0.975350813525687
  $ws.Cells.Item(3, 11).Value2 = 3
  $ws.Cells.Item(3, 12).Value2 = 1
  $ws.Cells.Item(3, 13).Value2 = 14.637306
  $ws.Cells.Item(3, 14).Value2 = 43.911918
  $ws.Cells.Item(3, 15).Value2 = 0.6566265559283671
  $ws.Cells.Item(3, 16).Value2 = 0.6566265559283672
  $ws.Cells.Item(3, 17).Value2 = 70.73273202216001
  $ws.Cells.Item(3, 18).Value2 = 636.5945881994401
  $ws.Cells.Item(3, 19).Value2 = 0.6404412455073027
  $ws.Cells.Item(3, 20).Value2 = 0.6404412455073029

# Row 4
  $ws.Cells.Item(4, 1).Value2 = "FAPs"
  $ws.Cells.Item(4, 2).Value2 = "Rspo3"
  $ws.Cells.Item(4, 3).Value2 = "Lgr4"
  $ws.Cells.Item(4, 4).Value2 = "M2"
  $ws.Cells.Item(4, 5).Value2 = 3
  $ws.Cells.Item(4, 6).Value2 = 1
  $ws.Cells.Item(4, 7).Value2 = 4.83236
  $ws.Cells.Item(4, 8).Value2 = 14.49708
  $ws.Cells.Item(4, 9).Value2 = 0.975350813525687
  $ws.Cells.Item(4, 10).Value2 = 0.975350813525687
  $ws.Cells.Item(4, 11).Value2 = 1
  $ws.Cells.Item(4, 12).Value2 = 0.3333333333333333
  $ws.Cells.Item(4, 13).Value2 = 0.02548533333333333
  $ws.Cells.Item(4, 14).Value2 = 0.076456
  $ws.Cells.Item(4, 15).Value2 = 0.001143266845234572
  $ws.Cells.Item(4, 16).Value2 = 0.001143266845234572
  $ws.Cells.Item(4, 17).Value2 = 0.1231543053866667
  $ws.Cells.Item(4, 18).Value2 = 1.10838874848
  $ws.Cells.Item(4, 19).Value2 = 0.001115086247576486
  $ws.Cells.Item(4, 20).Value2 = 0.001115086247576486

# Row 5
  $ws.Cells.Item(5, 1).Value2 = "FAPs"
  $ws.Cells.Item(5, 2).Value2 = "Rspo3"
  $ws.Cells.Item(5, 3).Value2 = "Lgr4"
  $ws.Cells.Item(5, 4).Value2 = "sCs"
  $ws.Cells.Item(5, 5).Value2 = 3
  $ws.Cells.Item(5, 6).Value2 = 1
  $ws.Cells.Item(5, 7).Value2 = 4.83236
  $ws.Cells.Item(5, 8).Value2 = 14.49708
  $ws.Cells.Item(5, 9).Value2 = 0.975350813525687
  $ws.Cells.Item(5, 10).Value2 = 0.975350813525687
  $ws.Cells.Item(5, 11).Value2 = 3
  $ws.Cells.Item(5, 12).Value2 = 1
  $ws.Cells.Item(5, 13).Value2 = 5.177582333333334
  $ws.Cells.Item(5, 14).Value2 = 15.532747
  $ws.Cells.Item(5, 15).Value2 = 0.2322652853996647
  $ws.Cells.Item(5, 16).Value2 = 0.2322652853996648
  $ws.Cells.Item(5, 17).Value2 = 25.01994176430667
  $ws.Cells.Item(5, 18).Value2 = 225.17947587876
  $ws.Cells.Item(5, 19).Value2 = 0.2265401350683389
  $ws.Cells.Item(5, 20).Value2 = 0.2265401350683389

# Row 6
  $ws.Cells.Item(6, 1).Value2 = "sCs"
  $ws.Cells.Item(6, 2).Value2 = "Rspo3"
  $ws.Cells.Item(6, 3).Value2 = "Lgr4"
  $ws.Cells.Item(6, 4).Value2 = "ECs"
  $ws.Cells.Item(6, 5).Value2 = 2
  $ws.Cells.Item(6, 6).Value2 = 0.6666666666666666
  $ws.Cells.Item(6, 7).Value2 = 0.122124
  $ws.Cells.Item(6, 8).Value2 = 0.366372
  $ws.Cells.Item(6, 9).Value2 = 0.02464918647431296
  $ws.Cells.Item(6, 10).Value2 = 0.02464918647431296
  $ws.Cells.Item(6, 11).Value2 = 3
  $ws.Cells.Item(6, 12).Value2 = 1
  $ws.Cells.Item(6, 13).Value2 = 2.451301666666667
  $ws.Cells.Item(6, 14).Value2 = 7.353905
  $ws.Cells.Item(6, 15).Value2 = 0.1099648918267336
  $ws.Cells.Item(6, 16).Value2 = 0.1099648918267337
  $ws.Cells.Item(6, 17).Value2 = 0.29936276474
  $ws.Cells.Item(6, 18).Value2 = 2.69426488266
  $ws.Cells.Item(6, 19).Value2 = 0.002710545124264811
  $ws.Cells.Item(6, 20).Value2 = 0.002710545124264812

# Row 7
  $ws.Cells.Item(7, 1).Value2 = "sCs"
  $ws.Cells.Item(7, 2).Value2 = "Rspo3"
  $ws.Cells.Item(7, 3).Value2 = "Lgr4"
  $ws.Cells.Item(7, 4).Value2 = "FAPs"
  $ws.Cells.Item(7, 5).Value2 = 2
  $ws.Cells.Item(7, 6).Value2 = 0.6666666666666666
  $ws.Cells.Item(7, 7).Value2 = 0.122124
  $ws.Cells.Item(7, 8).Value2 = 0.366372
  $ws.Cells.Item(7, 9).Value2 = 0.02464918647431296
  $ws.Cells.Item(7, 10).Value2 = 0.02464918647431296
  $ws.Cells.Item(7, 11).Value2 = 3
  $ws.Cells.Item(7, 12).Value2 = 1
  $ws.Cells.Item(7, 13).Value2 = 14.637306
  $ws.Cells.Item(7, 14).Value2 = 43.911918
  $ws.Cells.Item(7, 15).Value2 = 0.6566265559283671
  $ws.Cells.Item(7, 16).Value2 = 0.6566265559283672
  $ws.Cells.Item(7, 17).Value2 = 1.787566357944
  $ws.Cells.Item(7, 18).Value2 = 16.088097221496
  $ws.Cells.Item(7, 19).Value2 = 0.01618531042106421
  $ws.Cells.Item(7, 20).Value2 = 0.01618531042106421

# Row 8
  $ws.Cells.Item(8, 1).Value2 = "sCs"
  $ws.Cells.Item(8, 2).Value2 = "Rspo3"
  $ws.Cells.Item(8, 3).Value2 = "Lgr4"
  $ws.Cells.Item(8, 4).Value2 = "M2"
  $ws.Cells.Item(8, 5).Value2 = 2
  $ws.Cells.Item(8, 6).Value2 = 0.6666666666666666
  $ws.Cells.Item(8, 7).Value2 = 0.122124
  $ws.Cells.Item(8, 8).Value2 = 0.366372
  $ws.Cells.Item(8, 9).Value2 = 0.02464918647431296
  $ws.Cells.Item(8, 10).Value2 = 0.02464918647431296
  $ws.Cells.Item(8, 11).Value2 = 1
  $ws.Cells.Item(8, 12).Value2 = 0.3333333333333333
  $ws.Cells.Item(8, 13).Value2 = 0.02548533333333333
  $ws.Cells.Item(8, 14).Value2 = 0.076456
  $ws.Cells.Item(8, 15).Value2 = 0.001143266845234572
  $ws.Cells.Item(8, 16).Value2 = 0.001143266845234572
  $ws.Cells.Item(8, 17).Value2 = 0.003112370848
  $ws.Cells.Item(8, 18).Value2 = 0.028011337632
  $ws.Cells.Item(8, 19).Value2 = 0.0000281805976580864708347
  $ws.Cells.Item(8, 20).Value2 = 0.0000281805976580864708347

# Row 9
  $ws.Cells.Item(9, 1).Value2 = "sCs"
  $ws.Cells.Item(9, 2).Value2 = "Rspo3"
  $ws.Cells.Item(9, 3).Value2 = "Lgr4"
  $ws.Cells.Item(9, 4).Value2 = "sCs"
  $ws.Cells.Item(9, 5).Value2 = 2
  $ws.Cells.Item(9, 6).Value2 = 0.6666666666666666
  $ws.Cells.Item(9, 7).Value2 = 0.122124
  $ws.Cells.Item(9, 8).Value2 = 0.366372
  $ws.Cells.Item(9, 9).Value2 = 0.02464918647431296
  $ws.Cells.Item(9, 10).Value2 = 0.02464918647431296
  $ws.Cells.Item(9, 11).Value2 = 3
  $ws.Cells.Item(9, 12).Value2 = 1
  $ws.Cells.Item(9, 13).Value2 = 5.177582333333334
  $ws.Cells.Item(9, 14).Value2 = 15.532747
  $ws.Cells.Item(9, 15).Value2 = 0.2322652853996647
  $ws.Cells.Item(9, 16).Value2 = 0.2322652853996648
  $ws.Cells.Item(9, 17).Value2 = 0.632307064876
  $ws.Cells.Item(9, 18).Value2 = 5.690763583883999
  $ws.Cells.Item(9, 19).Value2 = 0.005725150331325856
  $ws.Cells.Item(9, 20).Value2 = 0.005725150331325857

